$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "jersson romero"
$ws.Range("C5").Select()
